# [MOSIP-21520] updated mosip-data for hindi, kannada, and tamil language
# Trim the stray leading space from the registration-center "name" / "addr_line1"
# shared-string values for the kan/hin/tam language rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C36").Value = 'ವರ್ಚುವಲ್ ಸೆಂಟರ್'
$ws.Range("C37").Value = 'ಸೆಂಟರ್ ಎ ಬೆನ್ ಮನ್ಸೂರ್'
$ws.Range("C38").Value = 'ಗ್ರಾಮೀಣ ಪುರಸಭೆ ಮ್ನಾಸ್ರಾ'
$ws.Range("C40").Value = 'ಕೇಂದ್ರ ಅಸ್ಸಾಂ'
$ws.Range("C41").Value = 'ಸೆಂಟರ್ ಮೆಹದಿಯಾ'
$ws.Range("C42").Value = 'ಸೆಂಟರ್ Ouled Oujih'
$ws.Range("C44").Value = 'ಸೆಂಟರ್ ಸಿಡಿ ಅಲ್ಲಲ್ ತಾಜಿ'
$ws.Range("C45").Value = 'ಕೇಂದ್ರ ಅಗ್ಡಾಲ್'
$ws.Range("C46").Value = 'ಕೇಂದ್ರ ಹಾಸನ'
$ws.Range("C47").Value = 'ಸೆಂಟರ್ ಸೌಸಿ'
$ws.Range("C49").Value = 'ಸೆಂಟರ್ ಹೇ ರಿಯಾಡ್'
$ws.Range("C50").Value = 'ಸೆಂಟರ್ ಮದೀನಾ'
$ws.Range("C51").Value = 'ಸೆಂಟರ್ ಯೂಸೌಫಿಯಾ'
$ws.Range("C52").Value = 'ಕೇಂದ್ರ ರಬತ್'
$ws.Range("C53").Value = 'आभासी केंद्र'
$ws.Range("C54").Value = 'सेंटर ए बेन मंसूर'
$ws.Range("C55").Value = 'ग्रामीण नगर मनसराय'
$ws.Range("C58").Value = 'केंद्र मेहदिया'
$ws.Range("C59").Value = 'केंद्र औलेड औजिहो'
$ws.Range("C60").Value = 'केंद्र सिदी ताइबिक'
$ws.Range("C61").Value = 'केंद्र सिदी अल्लाल ताज़िक'
$ws.Range("C62").Value = 'केंद्र Agdal'
$ws.Range("C63").Value = 'केंद्र हसन'
$ws.Range("C64").Value = 'केंद्र सूसी'
$ws.Range("C65").Value = 'केंद्र मदीनत अल इरफान'
$ws.Range("C67").Value = 'केंद्र मदीना'
$ws.Range("C68").Value = 'केंद्र युसूफिया'
$ws.Range("C69").Value = 'केंद्र रबातो'
$ws.Range("C74").Value = 'அசாம் மையம்'
$ws.Range("E74").Value = '7 கிமீ டேன்ஜியர் சாலை'
$ws.Range("C75").Value = 'மையம் மெஹதியா'
$ws.Range("C76").Value = 'மையம் Ouled Oujih'
$ws.Range("C77").Value = 'மையம் சிடி தைபி'
$ws.Range("C79").Value = 'மையம் அக்டல்'
$ws.Range("C80").Value = 'மையம் ஹாசன்'
$ws.Range("C81").Value = 'மையம் Souissi'
$ws.Range("C83").Value = 'மையம் ஹே ரியாட்'
$ws.Range("C84").Value = 'மதீனா மையம்'
$ws.Range("C85").Value = 'மையம் யூசுஃபியா'
$ws.Range("C86").Value = 'மையம் ரபாத்'

# Restore the view/selection state recorded for the sheet (scrolled to
# A66, active cell E86) instead of the prior topLeftCell=E71 / F86.
$win = $excel.ActiveWindow
$win.ScrollRow = 66
$win.ScrollColumn = 1
$ws.Range("E86").Select()

